# Updated testing data for 2020-08-03
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column map: A=ISO code B=Entity C=Date D=Source URL E=Source label F=Notes
#             G=Number of observations H=Cumulative total I=Cumulative total per thousand
#             J=Daily change in cumulative total K=Daily change per thousand
#             L=7-day smoothed daily change M=7-day smoothed per thousand

# --- Row 14: Chile - tests performed ---
$ws.Range("C14").Value = 44046
$ws.Range("G14").Value = 125
$ws.Range("H14").Value = 1697558
$ws.Range("I14").Value = 88.802
$ws.Range("J14").Value = 24269
$ws.Range("K14").Value = 1.27
$ws.Range("L14").Value = 21637
$ws.Range("M14").Value = 1.132

# --- Row 16: Costa Rica - people tested ---
$ws.Range("C16").Value = 44045
$ws.Range("G16").Value = 145
$ws.Range("H16").Value = 80406
$ws.Range("I16").Value = 15.784
$ws.Range("J16").Value = 2015
$ws.Range("K16").Value = 0.396
$ws.Range("L16").Value = 1731
$ws.Range("M16").Value = 0.34

# --- Row 17: Cote d'Ivoire - samples tested ---
$ws.Range("D17").Value = "https://www.facebook.com/Mshpci/posts/1658669570965397"
$ws.Range("C17").Value = 44045
$ws.Range("G17").Value = 110
$ws.Range("H17").Value = 102685
$ws.Range("I17").Value = 3.893
$ws.Range("J17").Value = 1102
$ws.Range("K17").Value = 0.042
$ws.Range("L17").Value = 1083
$ws.Range("M17").Value = 0.041

# --- Row 20: Czech Republic - tests performed ---
$ws.Range("C20").Value = 44045
$ws.Range("G20").Value = 183
$ws.Range("H20").Value = 706461
$ws.Range("I20").Value = 65.969
$ws.Range("J20").Value = 2571
$ws.Range("K20").Value = 0.24
$ws.Range("L20").Value = 6352
$ws.Range("M20").Value = 0.593

# --- Row 29: France - people tested ---
$ws.Range("C29").Value = 44043
$ws.Range("G29").Value = 80
$ws.Range("J29").Value = 87328
$ws.Range("K29").Value = 1.338
$ws.Range("L29").Value = 72890
$ws.Range("M29").Value = 1.117

# --- Row 32: Ghana - units unclear ---
$ws.Range("C32").Value = 44043
$ws.Range("G32").Value = 99
$ws.Range("H32").Value = 399446
$ws.Range("I32").Value = 12.855
$ws.Range("J32").Value = 2715
$ws.Range("K32").Value = 0.087
$ws.Range("L32").Value = 3548
$ws.Range("M32").Value = 0.114

# --- Row 33: Greece - samples tested ---
$ws.Range("D33").Value = "https://eody.gov.gr/covid-gr-daily-report-20200803"
$ws.Range("C33").Value = 44046
$ws.Range("G33").Value = 133
$ws.Range("H33").Value = 564856
$ws.Range("I33").Value = 54.193
$ws.Range("J33").Value = 19417
$ws.Range("K33").Value = 1.863
$ws.Range("L33").Value = 16220
$ws.Range("M33").Value = 1.556

# --- Row 62: Nigeria - samples tested ---
$ws.Range("C62").Value = 44045
$ws.Range("G62").Value = 88
$ws.Range("H62").Value = 287532
$ws.Range("I62").Value = 1.395
$ws.Range("J62").Value = 1441
$ws.Range("K62").Value = 0.007
$ws.Range("L62").Value = 3565
$ws.Range("M62").Value = 0.017

# --- Row 70: Poland - people tested ---
$ws.Range("D70").Value = "https://twitter.com/MZ_GOV_PL/status/1290210801614110723"
$ws.Range("C70").Value = 44046
$ws.Range("G70").Value = 97
$ws.Range("H70").Value = 1980641
$ws.Range("I70").Value = 52.333
$ws.Range("J70").Value = 11071
$ws.Range("K70").Value = 0.293
$ws.Range("L70").Value = 18997
$ws.Range("M70").Value = 0.502

# --- Row 71: Poland - samples tested ---
$ws.Range("D71").Value = "https://twitter.com/MZ_GOV_PL/status/1290210801614110723"
$ws.Range("C71").Value = 44046
$ws.Range("G71").Value = 148
$ws.Range("H71").Value = 2315210
$ws.Range("I71").Value = 61.174
$ws.Range("J71").Value = 17783
$ws.Range("K71").Value = 0.47
$ws.Range("L71").Value = 25839
$ws.Range("M71").Value = 0.683

# --- Row 73: Qatar - people tested ---
$ws.Range("C73").Value = 44046
$ws.Range("G73").Value = 136
$ws.Range("H73").Value = 502792
$ws.Range("I73").Value = 174.516
$ws.Range("J73").Value = 2256
$ws.Range("K73").Value = 0.783
$ws.Range("L73").Value = 3657
$ws.Range("M73").Value = 1.269

# --- Row 77: Saudi Arabia - units unclear ---
$ws.Range("C77").Value = 44046
$ws.Range("G77").Value = 88
$ws.Range("H77").Value = 3464427
$ws.Range("I77").Value = 99.513
$ws.Range("J77").Value = 41361
$ws.Range("K77").Value = 1.188
$ws.Range("L77").Value = 51852
$ws.Range("M77").Value = 1.489

# --- Row 86: Spain - tests performed ---
$ws.Range("D86").Value = "https://www.mscbs.gob.es/profesionales/saludPublica/ccayes/alertasActual/nCov-China/documentos/COVID-19_pruebas_diagnosticas_30_07_2020.pdf"
$ws.Range("C86").Value = 44042
$ws.Range("G86").Value = 16
$ws.Range("H86").Value = 4652493
$ws.Range("I86").Value = 99.508
$ws.Range("L86").Value = 43639
$ws.Range("M86").Value = 0.933

# --- Row 90: Taiwan - people tested ---
$ws.Range("C90").Value = 44046
$ws.Range("G90").Value = 200
$ws.Range("H90").Value = 82337
$ws.Range("I90").Value = 3.457
$ws.Range("J90").Value = 273
$ws.Range("K90").Value = 0.011
$ws.Range("L90").Value = 216

# --- Row 91: Thailand - people tested ---
$ws.Range("D91").Value = "https://ddc.moph.go.th/viralpneumonia/file/situation/situation-no213-030863.pdf"
$ws.Range("C91").Value = 44046
$ws.Range("G91").Value = 141
$ws.Range("H91").Value = 379104
$ws.Range("I91").Value = 5.431
$ws.Range("J91").Value = 1424
$ws.Range("K91").Value = 0.02
$ws.Range("L91").Value = 1795
$ws.Range("M91").Value = 0.026

# --- Row 92: Thailand - tests performed ---
$ws.Range("D92").Value = "https://ddc.moph.go.th/viralpneumonia/file/situation/situation-no213-030863.pdf"
$ws.Range("C92").Value = 44046
$ws.Range("G92").Value = 46
$ws.Range("H92").Value = 731449
$ws.Range("I92").Value = 10.479
$ws.Range("J92").Value = 1424
$ws.Range("K92").Value = 0.02
$ws.Range("L92").Value = 4302
$ws.Range("M92").Value = 0.062

# --- Row 94: Tunisia - tests performed ---
$ws.Range("C94").Value = 44044
$ws.Range("G94").Value = 144
$ws.Range("H94").Value = 93209
$ws.Range("I94").Value = 7.887
$ws.Range("J94").Value = 249
$ws.Range("K94").Value = 0.021
$ws.Range("L94").Value = 909
$ws.Range("M94").Value = 0.077

# --- Row 95: Turkey - tests performed ---
$ws.Range("C95").Value = 44045
$ws.Range("G95").Value = 128
$ws.Range("H95").Value = 4927217
$ws.Range("I95").Value = 58.422
$ws.Range("J95").Value = 41301
$ws.Range("K95").Value = 0.49
$ws.Range("L95").Value = 50647

# --- Row 96: Uganda - samples tested ---
$ws.Range("D96").Value = "https://twitter.com/MinofHealthUG/status/1290255151530418176/photo/2"
$ws.Range("C96").Value = 44045
$ws.Range("G96").Value = 32
$ws.Range("H96").Value = 280747
$ws.Range("I96").Value = 6.138
$ws.Range("J96").Value = 2485
$ws.Range("K96").Value = 0.054
$ws.Range("L96").Value = 2639
$ws.Range("M96").Value = 0.058

# --- Row 98: United Arab Emirates - samples tested ---
$ws.Range("C98").Value = 44046
$ws.Range("G98").Value = 157
$ws.Range("H98").Value = 5189891
$ws.Range("I98").Value = 524.74
$ws.Range("J98").Value = 27811
$ws.Range("K98").Value = 2.812
$ws.Range("L98").Value = 44942
$ws.Range("M98").Value = 4.544

# --- Row 99: United Kingdom - tests performed ---
$ws.Range("D99").Value = "https://assets.publishing.service.gov.uk/government/uploads/system/uploads/attachment_data/file/906445/2020-08-03-COVID-19-UK-testing-time-series.csv"
$ws.Range("C99").Value = 44045
$ws.Range("G99").Value = 125
$ws.Range("H99").Value = 9818696
$ws.Range("I99").Value = 144.635
$ws.Range("J99").Value = 124088
$ws.Range("K99").Value = 1.828
$ws.Range("L99").Value = 130584
$ws.Range("M99").Value = 1.924

$wb.Save()
